# Update the acquisition timestamp (column A) for all data rows on the
# "ランサーズ" sheet from "2025-10-29 12:50:21" to "2025-10-29 18:28:38".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldValue = "2025-10-29 12:50:21"
$newValue = "2025-10-29 18:28:38"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
